# Generate Report for Handoff
#
# A fresh handoff run produced a new working-file GUID and new content
# hashes for the generated .xlf payloads; refresh the localization-status
# report (Overview + per-locale sheets) with the new identifiers and
# timestamps, keeping each hyperlink's target intact but refreshing its
# displayed text.

$wb = $excel.ActiveWorkbook

$oldGuid = "51fb5b97-ff4e-4918-ab33-1c47943dbc3c"
$newGuid = "8870f70d-9f1b-4ed2-befb-a238c465f730"

$newZhHash = "96903885617b4de9c623c4d6340d2b89eeb67058"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e86c7e8321f16caa6a1c3f8d2de013713ae5eac8/e2e/$oldGuid.md"

function Update-Hyperlink($range, $displayText) {
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $hyperlinkAddress, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
Update-Hyperlink $wsOverview.Range("B2") "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-29 21:12:26"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
Update-Hyperlink $wsZh.Range("A2") "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-29 21:12:22"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
Update-Hyperlink $wsDe.Range("A2") "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
